$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 2. Data reporter - updated contact details
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Update active cell selection to B9
$ws.Range("B9").Select()

# Update workbook window view/geometry settings (restored/maximized full view)
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 28800
$win.Height = 11835
